$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the club name typo in cell A4: "Athlético" -> "Athletico"
$ws.Range("A4").Value = "Athletico"

# Update the active selection to cell A5
$ws.Range("A5").Select()
